$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between rows 3/4 and rows 5/6 (per diff): A, B, E, F, G, H, Q, R, AC
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Read the current contents of rows 3-6 for the affected columns into memory first,
# so writes don't clobber values we still need to read.
$data = @{}
foreach ($r in 3..6) {
    $data[$r] = @{}
    foreach ($col in $cols) {
        $data[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# The edit swaps the full record contents of rows (3,4) with rows (5,6):
#   new row3 <- old row5
#   new row4 <- old row6
#   new row5 <- old row3
#   new row6 <- old row4
$mapping = @{ 3 = 5; 4 = 6; 5 = 3; 6 = 4 }

foreach ($destRow in 3..6) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $val = $data[$srcRow][$col]
        if ($val -eq $null) {
            $ws.Range("$col$destRow").Value = ""
        } else {
            $ws.Range("$col$destRow").Value = $val
        }
    }
}
